$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 254.41667
$ws.Range("I5").Value = 105.3
$ws.Range("K5").Value = 105.3
$ws.Range("M5").Value = 9.700000000000003
$ws.Range("H18").Value = 310.9524
$ws.Range("I18").Value = 306.5
$ws.Range("J18").Value = 400
$ws.Range("K18").Value = 306.5
$ws.Range("L18").Value = 400
$ws.Range("M18").Value = -22.5
$ws.Range("N18").Value = -968
$ws.Range("H28").Value = 571.3570999999999
$ws.Range("I28").Value = 207.72728
$ws.Range("J28").Value = 1904.6666
$ws.Range("K28").Value = 207.72728
$ws.Range("L28").Value = 1904.6666
$ws.Range("M28").Value = 277.27272
$ws.Range("N28").Value = -2874.6666
$ws.Range("H40").Value = 1945.8572
$ws.Range("I40").Value = 1826.5333
$ws.Range("J40").Value = 2083.5386
$ws.Range("K40").Value = 1826.5333
$ws.Range("L40").Value = 2083.5386
$ws.Range("M40").Value = -1651.5333
$ws.Range("N40").Value = -2433.5386
$ws.Range("H41").Value = 240.21428
$ws.Range("I41").Value = 270.57144
$ws.Range("J41").Value = 209.85715
$ws.Range("K41").Value = 270.57144
$ws.Range("L41").Value = 209.85715
$ws.Range("M41").Value = 169.42856
$ws.Range("N41").Value = -1089.85715
$ws.Range("H64").Value = 5166
$ws.Range("I64").Value = 3949.25
$ws.Range("K64").Value = 3949.25
$ws.Range("M64").Value = -3701.25
$ws.Range("H67").Value = 5166
$ws.Range("I67").Value = 3949.25
$ws.Range("K67").Value = 3949.25
$ws.Range("M67").Value = -3091.25
$ws.Range("H86").Value = 3085.8823
$ws.Range("I86").Value = 3160.2
$ws.Range("J86").Value = 3054.9167
$ws.Range("K86").Value = 3160.2
$ws.Range("L86").Value = 3054.9167
$ws.Range("M86").Value = -2037.2
$ws.Range("N86").Value = -5300.9167
$ws.Range("H89").Value = 3085.8823
$ws.Range("I89").Value = 3160.2
$ws.Range("J89").Value = 3054.9167
$ws.Range("K89").Value = 15801
$ws.Range("L89").Value = 15274.5835
$ws.Range("M89").Value = -10185
$ws.Range("N89").Value = -26506.5835
$ws.Range("H92").Value = 392.4516
$ws.Range("I92").Value = 364.625
$ws.Range("J92").Value = 487.85715
$ws.Range("K92").Value = 364.625
$ws.Range("L92").Value = 487.85715
$ws.Range("M92").Value = 883.375
$ws.Range("N92").Value = -2983.85715
$ws.Range("H112").Value = 2114.9092
$ws.Range("I112").Value = 957
$ws.Range("J112").Value = 2549.125
$ws.Range("K112").Value = 2871
$ws.Range("L112").Value = 7647.375
$ws.Range("M112").Value = -1763
$ws.Range("N112").Value = -9863.375
$ws.Range("H124").Value = 105000
$ws.Range("J124").Value = 105000
$ws.Range("L124").Value = 105000
$ws.Range("N124").Value = -114820
$ws.Range("H137").Value = 16885.732
$ws.Range("I137").Value = 7460.533
$ws.Range("J137").Value = 26310.934
$ws.Range("K137").Value = 22381.599
$ws.Range("L137").Value = 78932.802
$ws.Range("M137").Value = -19831.599
$ws.Range("N137").Value = -84032.802
$ws.Range("H139").Value = 125111
$ws.Range("J139").Value = 125111
$ws.Range("L139").Value = 125111
$ws.Range("N139").Value = -135391
$ws.Range("H141").Value = 2353.8
$ws.Range("I141").Value = 2353.8
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 7061.400000000001
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").Value = -1881.400000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 10349
$ws.Range("J45").Value = 2673.8333
$ws.Range("L45").Value = 2673.8333
$ws.Range("N45").Value = -3427.8333
$ws.Range("H74").Value = 16949.902
$ws.Range("I74").Value = 18194.75
$ws.Range("J74").Value = 5331.3335
$ws.Range("K74").Value = 18194.75
$ws.Range("L74").Value = 5331.3335
$ws.Range("M74").Value = -17320.75
$ws.Range("N74").Value = -7079.3335
$ws.Range("H77").Value = 16949.902
$ws.Range("I77").Value = 18194.75
$ws.Range("J77").Value = 5331.3335
$ws.Range("K77").Value = 90973.75
$ws.Range("L77").Value = 26656.6675
$ws.Range("M77").Value = -86605.75
$ws.Range("N77").Value = -35392.6675
$ws.Range("H112").Value = 37249.5
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -846
$ws.Range("H86").Value = 186912.8
$ws.Range("I86").Value = 323318.16
$ws.Range("J86").Value = 3062.087
$ws.Range("K86").Value = 323318.16
$ws.Range("L86").Value = 3062.087
$ws.Range("M86").Value = -322195.16
$ws.Range("N86").Value = -5308.087
$ws.Range("H89").Value = 186912.8
$ws.Range("I89").Value = 323318.16
$ws.Range("J89").Value = 3062.087
$ws.Range("K89").Value = 1616590.8
$ws.Range("L89").Value = 15310.435
$ws.Range("M89").Value = -1610974.8
$ws.Range("N89").Value = -26542.435

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 499998.5
$ws.Range("J17").Value = 499998.5
$ws.Range("L17").Value = 499998.5
$ws.Range("N17").Value = -500346.5
$ws.Range("H31").Value = 4082.5
$ws.Range("I31").Value = 1999
$ws.Range("J31").Value = 4499.2
$ws.Range("K31").Value = 1999
$ws.Range("L31").Value = 4499.2
$ws.Range("M31").Value = -1704
$ws.Range("N31").Value = -5089.2
$ws.Range("H34").Value = 4082.5
$ws.Range("I34").Value = 1999
$ws.Range("J34").Value = 4499.2
$ws.Range("K34").Value = 1999
$ws.Range("L34").Value = 4499.2
$ws.Range("M34").Value = -1797
$ws.Range("N34").Value = -4903.2
$ws.Range("H58").Value = 8844.182000000001
$ws.Range("I58").Value = 4673.375
$ws.Range("K58").Value = 4673.375
$ws.Range("M58").Value = -4470.375
$ws.Range("H132").Value = 30514.045
$ws.Range("I132").Value = 22323.363
$ws.Range("K132").Value = 66970.08900000001
$ws.Range("M132").Value = -64440.08900000001
$ws.Range("H136").Value = 8844.182000000001
$ws.Range("I136").Value = 4673.375
$ws.Range("K136").Value = 14020.125
$ws.Range("M136").Value = -11470.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2314.2083
$ws.Range("I5").Value = 1515.3334
$ws.Range("J5").Value = 2428.3333
$ws.Range("K5").Value = 4546.0002
$ws.Range("L5").Value = 7284.999899999999
$ws.Range("M5").Value = -4434.0002
$ws.Range("N5").Value = -7508.999899999999
$ws.Range("H122").Value = 7692832
$ws.Range("J122").Value = 25000778
$ws.Range("L122").Value = 225007002
$ws.Range("N122").Value = -225011902
$ws.Range("H135").Value = 2314.2083
$ws.Range("I135").Value = 1515.3334
$ws.Range("J135").Value = 2428.3333
$ws.Range("K135").Value = 13638.0006
$ws.Range("L135").Value = 21854.9997
$ws.Range("M135").Value = -11103.0006
$ws.Range("N135").Value = -26924.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 29899
$ws.Range("I20").Value = 29899
$ws.Range("K20").Value = 29899
$ws.Range("M20").Value = -29654
$ws.Range("H54").Value = 750
$ws.Range("J54").Value = 750
$ws.Range("L54").Value = 750
$ws.Range("N54").Value = -1530
$ws.Range("H95").Value = 18107.25
$ws.Range("J95").Value = 18107.25
$ws.Range("L95").Value = 18107.25
$ws.Range("N95").Value = -23599.25
$ws.Range("H113").Value = 120315
$ws.Range("I113").Value = 127665.94
$ws.Range("J113").Value = 2700
$ws.Range("K113").Value = 127665.94
$ws.Range("L113").Value = 2700
$ws.Range("M113").Value = -125495.94
$ws.Range("N113").Value = -7040
$ws.Range("H132").Value = 44224.24
$ws.Range("I132").Value = 24697.732
$ws.Range("J132").Value = 93040.5
$ws.Range("K132").Value = 74093.196
$ws.Range("L132").Value = 279121.5
$ws.Range("M132").Value = -71563.196
$ws.Range("N132").Value = -284181.5
$ws.Range("H135").Value = 74999
$ws.Range("J135").Value = 74999
$ws.Range("L135").Value = 74999
$ws.Range("N135").Value = -85139
$ws.Range("H140").Value = 59833
$ws.Range("J140").Value = 59799.6
$ws.Range("L140").Value = 59799.6
$ws.Range("N140").Value = -70159.60000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2744.9375
$ws.Range("I22").Value = 2391.5
$ws.Range("K22").Value = 2391.5
$ws.Range("M22").Value = -2096.5
$ws.Range("H27").Value = 2744.9375
$ws.Range("I27").Value = 2391.5
$ws.Range("K27").Value = 2391.5
$ws.Range("M27").Value = -2284.5
$ws.Range("H46").Value = 1885.64
$ws.Range("I46").Value = 1172.4546
$ws.Range("J46").Value = 2446
$ws.Range("K46").Value = 1172.4546
$ws.Range("L46").Value = 2446
$ws.Range("M46").Value = -984.4546
$ws.Range("N46").Value = -2822
$ws.Range("H132").Value = 5317.978
$ws.Range("I132").Value = 4870.054
$ws.Range("K132").Value = 14610.162
$ws.Range("M132").Value = -12080.162
$ws.Range("H136").Value = 5434.2764
$ws.Range("I136").Value = 4809.2573
$ws.Range("J136").Value = 7257.25
$ws.Range("K136").Value = 14427.7719
$ws.Range("L136").Value = 21771.75
$ws.Range("M136").Value = -11877.7719
$ws.Range("N136").Value = -26871.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 435.88
$ws.Range("I113").Value = 294.08694
$ws.Range("K113").Value = 882.2608200000001
$ws.Range("M113").Value = 1287.73918

